$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(22,3).Value = 5
$ws.Cells.Item(22,4).Value = 5
